$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.599.71'
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").Value = '1.815.35'
$ws.Range("E3").Value = '  +1.43%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.37'
$ws.Range("E5").Value = '  +1.11%  '
$ws.Range("E6").Value = '  +0.97%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '34.78'
$ws.Range("E8").Value = '  +7.76%  '
$ws.Range("E9").Value = '  +1.97%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0695'
$ws.Range("E10").Value = '  +0.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0953'
$ws.Range("E11").Value = '  +0.44%  '
$ws.Range("D12").Value = '2.078.66'
$ws.Range("E12").Value = '  +1.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.31'
$ws.Range("E13").Value = '  +2.87%  '
$ws.Range("D14").Value = '1.832.04'
$ws.Range("E14").Value = '  +2.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.643'
$ws.Range("E15").Value = '  +2.58%  '
$ws.Range("D16").Value = '34.646.26'
$ws.Range("E16").Value = '  +1.27%  '
$ws.Range("E17").Value = '  +3.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.14'
$ws.Range("E18").Value = '  +1.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '247.36'
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("D20").Value = '0.0₃0802'
$ws.Range("E20").Value = '  -0.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.56'
$ws.Range("E21").Value = '  +5.65%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("E23").Value = '  +1.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '172.63'
$ws.Range("E24").Value = '  +6.33%  '
$ws.Range("E25").Value = '  +2.23%  '
$ws.Range("E26").Value = '  +4.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.77'
$ws.Range("E27").Value = '  +2.58%  '
$ws.Range("E28").Value = '  +1.58%  '
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.08'
$ws.Range("E30").Value = '  +6.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0532'
$ws.Range("E31").Value = '  +2.28%  '
$ws.Range("E32").Value = '  +2.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.24'
$ws.Range("E33").Value = '  +1.10%  '
$ws.Range("E34").Value = '  +2.53%  '
$ws.Range("E35").Value = '  -0.39%  '
$ws.Range("D36").Value = '1.418.39'
$ws.Range("E36").Value = '  -1.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.675'
$ws.Range("E37").Value = '  +1.91%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.06'
$ws.Range("E38").Value = '  +1.14%  '
$ws.Range("E39").Value = '  +1.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '85.97'
$ws.Range("E40").Value = '  +4.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.85'
$ws.Range("E41").Value = '  +4.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.961'
$ws.Range("E42").Value = '  +4.13%  '
$ws.Range("E43").Value = '  +0.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.88'
$ws.Range("E44").Value = '  -1.58%  '
$ws.Range("E46").Value = '  +2.86%  '
$ws.Range("E47").Value = '  +0.86%  '
$ws.Range("D48").Value = '1.978.75'
$ws.Range("E48").Value = '  +1.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.97'
$ws.Range("E49").Value = '  +0.40%  '
$ws.Range("D50").Value = '0.0₆0131'
$ws.Range("E50").Value = '  +1.40%  '
$ws.Range("E51").Value = '  +0.09%  '
